$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 2; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 3; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 4; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 29; I = 'ba'; J = 'Appreciation' },
    @{ Row = 30; I = 'ba'; J = 'Appreciation' },
    @{ Row = 33; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 42; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 49; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 54; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 59; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 66; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 67; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 83; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 124; I = 'ba'; J = 'Appreciation' },
    @{ Row = 134; I = '%'; J = 'Uninterpretable' },
    @{ Row = 146; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 150; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 152; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 162; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 177; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 189; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 195; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 196; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 240; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 259; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 267; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 273; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 275; I = '%'; J = 'Uninterpretable' },
    @{ Row = 277; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 286; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 287; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 311; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 319; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 325; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 352; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 354; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 356; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 358; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 360; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 379; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 387; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 392; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 393; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 423; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 427; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 438; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 439; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 443; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 446; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 454; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 462; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 483; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 492; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 510; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 515; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 517; I = 'ba'; J = 'Appreciation' },
    @{ Row = 520; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 523; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 524; I = 'ba'; J = 'Appreciation' },
    @{ Row = 526; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 529; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 536; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 537; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 550; I = 'sv'; J = 'Statement-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
